# Cursed-coin / cursed-pot style trap card effect tweak:
# the "Boulder" (滚石) card's death-trigger effect now only decrements the
# rank of the *other* cards sharing its destination slot (adds "其他"/"other"),
# so the moved card itself stays around longer instead of also losing a
# point. This lets it, and the slot it lands in, survive longer.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "死亡时：如果本牌点数大于1，则将本牌移动到相邻槽位中而不是送墓，然后使本牌和那个槽位中所有其他牌点数减1。<br>"

# Restore the sheet's scroll position / selection so the view lands back
# near the top of the table (selection on D9) instead of the previously
# scrolled-down D14 state.
$ws.Range("D9").Select()
